# Update loading_percent values on Sheet1 for rows 2-25.
# Columns updated: B, D, E, F, G, I, K, M, N (column indices 2,4,5,6,7,9,11,13,14)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column letters (in order) mapped to their 1-based column index
$cols = @(2, 4, 5, 6, 7, 9, 11, 13, 14)

$data = @(
    ,(2, 4.990939164366616, 6.507840982241279, 9.618012127186548, 35.52019024459483, 3.673815597944844, 22.04732753250851, 15.51992752456134, 16.55245628148672, 20.09511346026386)
    ,(3, 4.926753005089897, 6.531579384626017, 9.493302215870655, 35.22932169372265, 3.677683319687123, 22.09362646585049, 15.13835380996798, 16.31641927516037, 20.15031718537049)
    ,(4, 4.887998066698472, 6.546842910673605, 9.419053196848148, 35.06039453937645, 3.680179439977323, 22.12540637956185, 14.90425849524191, 16.17496664291349, 20.18612116220475)
    ,(5, 4.872390437883902, 6.553236201918351, 9.389414507076193, 34.99404489475332, 3.681227258455599, 22.13919803528738, 14.80906765839697, 16.11826912038409, 20.20119178115052)
    ,(6, 4.869810589602809, 6.554308276734658, 9.384531342483179, 34.98317955709946, 3.681403101405924, 22.14153887340739, 14.79327832834815, 16.10891370642628, 20.20372325728087)
    ,(7, 4.887786799501153, 6.546928430911449, 9.418650931260622, 35.05948957157163, 3.680193447048997, 22.125588975714, 14.90297367865235, 16.17419807847738, 20.18632246522423)
    ,(8, 4.968684036764326, 6.515883336391759, 9.57455087142322, 35.41792908028361, 3.675124080969409, 22.06259456072215, 15.38841335169359, 16.47039309890616, 20.11375171437296)
    ,(9, 5.131619422700456, 6.460451148125277, 9.897139049378055, 36.19471375917445, 3.666140174219825, 21.96574255728804, 16.335282790329, 17.07557060641579, 19.98657292281171)
    ,(10, 5.252735025288806, 6.423028562558724, 10.14226578494014, 36.80637305569483, 3.660115409898228, 21.91096067023519, 17.01948969477802, 17.53064376098246, 19.90234213983214)
    ,(11, 5.307886217678909, 6.406717944338693, 10.25507475628799, 37.09257839567319, 3.657497937132963, 21.88961737745928, 17.32665554504886, 17.73900667693052, 19.866019469106)
    ,(12, 5.3287578994862, 6.400643878552204, 10.29794025881503, 37.20201756042348, 3.65652435948944, 21.88205138656825, 17.44226209713271, 17.81802572647701, 19.85255168854276)
    ,(13, 5.324263707112404, 6.401947484420722, 10.28870248605448, 37.17840213434048, 3.656733255835167, 21.88365786902492, 17.41739762084493, 17.80100358846047, 19.85543945830361)
    ,(14, 5.309603718674265, 6.406216175507982, 10.25859860258891, 37.10156128320993, 3.65741748820254, 21.88898456208398, 17.33618147800935, 17.74550572237738, 19.8649057179018)
    ,(15, 5.300621767880027, 6.408844205424034, 10.24017710393017, 37.05462945360589, 3.657838889186707, 21.89231460004559, 17.28633821857153, 17.71152459019232, 19.87074143035879)
    ,(16, 5.249130082931377, 6.424108837572849, 10.13491620205041, 36.78782178908616, 3.660288937464626, 21.91242764755583, 16.99932276323483, 17.51704723082525, 19.90475604588218)
    ,(17, 5.217540739466532, 6.433655738090763, 10.07064654072925, 36.62612248137034, 3.661823439820724, 21.92568390780921, 16.82211155116376, 17.39803137143946, 19.92613372841228)
    ,(18, 5.199377782955716, 6.439213990243934, 10.03380501258819, 36.53387164519813, 3.662717650263092, 21.93364516104214, 16.71980679712466, 17.32970799635077, 19.93861731766619)
    ,(19, 5.193229859784928, 6.441107452431488, 10.02135373434263, 36.5027692091921, 3.663022411332838, 21.93639845907147, 16.68510723048114, 17.30659982603831, 19.94287628530616)
    ,(20, 5.220902969277054, 6.432632507812203, 10.07747554076727, 36.64325818783576, 3.661658889171369, 21.92423790557679, 16.84101598534647, 17.41068777890519, 19.9238386099627)
    ,(21, 5.313910225038864, 6.404959579004835, 10.2674371619748, 37.12410323872058, 3.657216035772941, 21.88740595794807, 17.36005684671801, 17.76180423511819, 19.86211746308765)
    ,(22, 5.374611561928782, 6.38747060396077, 10.39242801228885, 37.44449768286756, 3.654414921015546, 21.8663440623611, 17.69508819278815, 17.99192392629193, 19.82345101369779)
    ,(23, 5.342228545428953, 6.396750222972574, 10.32565415986452, 37.27296439459907, 3.655900584148888, 21.87730919381875, 17.51669753372516, 17.86907075907694, 19.84393502457981)
    ,(24, 5.219382909271232, 6.433094893385855, 10.07438781057485, 36.63550891002148, 3.661733245119145, 21.92489058413872, 16.8324706008241, 17.40496549839742, 19.92487563100913)
    ,(25, 5.08720693419875, 6.474865716814605, 9.808285333594032, 35.977069062351, 3.668468898349057, 21.98907568383997, 16.08056638073072, 16.90970062450446, 20.01936007264079)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Cells.Item($row, $col).Value = $value
    }
}
